$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 51299.77620886141
$ws.Range("C2").Value = 32470.46482413697
$ws.Range("D2").Value = 30562.23919850059
$ws.Range("B3").Value = 32633.97734808491
$ws.Range("C3").Value = 22890.15223737805
$ws.Range("D3").Value = 24229.09556265955
$ws.Range("B5").Value = 2314.123146287206
$ws.Range("C5").Value = -5269.701964419652
$ws.Range("D5").Value = -3930.758639138152
$ws.Range("B6").Value = 3596.202853792583
$ws.Range("C6").Value = 603.2279645160052
$ws.Range("D6").Value = 854.6824965449755
$ws.Range("B8").Value = 2028.702853792583
$ws.Range("C8").Value = -964.2720354839948
$ws.Range("D8").Value = -712.8175034550245
$ws.Range("B9").Value = -14.2883361816406
$ws.Range("C9").Value = -1.167364501953102
$ws.Range("D9").Value = -0.9222778320312273
$ws.Range("B10").Value = 26.99999999999983
$ws.Range("C10").Value = 23.72113514408119
$ws.Range("D10").Value = 26.79948907117119
$ws.Range("B11").Value = 26.46706608240203
$ws.Range("C11").Value = 23.99999999999983
$ws.Range("D11").Value = 27.0000000000004
$ws.Range("B12").Value = 26.88362262949369
$ws.Range("C12").Value = 23.25952019118625
$ws.Range("D12").Value = 26.34810447416294
$ws.Range("B13").Value = 64.73349270889196
$ws.Range("C13").Value = 22.31024172486366
$ws.Range("D13").Value = 25.5907792209062
$ws.Range("B14").Value = 11.34512000809707
$ws.Range("C14").Value = 16.04738618027744
$ws.Range("D14").Value = 16.51352254903583
$ws.Range("B15").Value = 8.423013137011665
$ws.Range("C15").Value = 6.578993866641933
$ws.Range("D15").Value = 6.397913401514893
$ws.Range("B16").Value = 15.6007048102111
$ws.Range("C16").Value = 12.46304758374061
$ws.Range("D16").Value = 12.34675695646267
$ws.Range("B17").Value = 15.19575318482976
$ws.Range("C17").Value = 12.99920415303144
$ws.Range("D17").Value = 12.80326246143624
$ws.Range("B18").Value = 1.214233069637331
$ws.Range("C18").Value = 7.163226421827404
$ws.Range("D18").Value = 6.869636926063983
$ws.Range("B19").Value = 1.999999999999918
$ws.Range("C19").Value = 2.120000000000013
$ws.Range("D19").Value = 1.760000000000054
$ws.Range("B20").Value = 0.2189051292047605
$ws.Range("C20").Value = 0.2446962585405372
$ws.Range("D20").Value = 0.2061140735987233
$ws.Range("B21").Value = 1.816094870795158
$ws.Range("C21").Value = 1.892303741459475
$ws.Range("D21").Value = 1.570885926401331
$ws.Range("B22").Value = 1.819285348705744
$ws.Range("C22").Value = 1.89535013941786
$ws.Range("D22").Value = 1.57388701945878
$ws.Range("B23").Value = 0.2189636826515198
$ws.Range("C23").Value = 0.2447559684514999
$ws.Range("D23").Value = 0.2061733305454254
$ws.Range("B25").Value = 1.038249015808105
$ws.Range("B26").Value = 0.1229048114135727
$ws.Range("C26").Value = 0.1203275423241102
$ws.Range("D26").Value = 0.1203868915115663
$ws.Range("B29").Value = 104.5951640625
$ws.Range("C29").Value = 79.60096875000001
$ws.Range("D29").Value = 78.5914296875
$ws.Range("B31").Value = 104.5951676965513
$ws.Range("C31").Value = 79.6009678006568
$ws.Range("D31").Value = 78.59142710051891
